$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","V","W","X","Y")

# Shift all data rows 11..114 down by one (into rows 12..115), working from the
# bottom up so we never overwrite a row before reading it.
for ($r = 115; $r -ge 12; $r--) {
    $src = $r - 1
    foreach ($col in $cols) {
        $srcCell = $ws.Range($col + $src)
        $dstCell = $ws.Range($col + $r)
        $dstCell.Formula = $srcCell.Formula
    }
}

# Write the brand-new record into row 11.
$ws.Range("A11").Formula = "A 40983-2023"
$ws.Range("B11").Formula = "45173"
$ws.Range("C11").Formula = "45202"
$ws.Range("D11").Formula = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E11").Formula = "ALE"
$ws.Range("F11").Formula = "Kyrkan"
$ws.Range("G11").Formula = "12.6"
$ws.Range("H11").Formula = "2"
$ws.Range("I11").Formula = "0"
$ws.Range("J11").Formula = "2"
$ws.Range("K11").Formula = "0"
$ws.Range("L11").Formula = "0"
$ws.Range("M11").Formula = "0"
$ws.Range("N11").Formula = "0"
$ws.Range("O11").Formula = "2"
$ws.Range("P11").Formula = "0"
$ws.Range("Q11").Formula = "2"
$ws.Range("R11").Formula = "Spillkråka`r`nTalltita"
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/artfynd/A 40983-2023.xlsx", "A 40983-2023")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/kartor/A 40983-2023.png", "A 40983-2023")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/klagomål/A 40983-2023.docx", "A 40983-2023")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/klagomålsmail/A 40983-2023.docx", "A 40983-2023")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/tillsyn/A 40983-2023.docx", "A 40983-2023")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/tillsynsmail/A 40983-2023.docx", "A 40983-2023")'

# Every record's "Förändrad" (last-changed) date is refreshed to 45202.
for ($r = 2; $r -le 115; $r++) {
    $ws.Range("C" + $r).Formula = "45202"
}
